$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# The paragraph currently holds two runs: the **ID__...__ID** placeholder
# and a trailing " " run. Delete the trailing space run entirely (rather
# than just clearing its text) by removing its Range. The paragraph mark
# occupies the last character slot (End-1), so the trailing space run
# sits in the slot just before that.
$full = $p1.Range
$spaceRange = $d.Range($full.End - 2, $full.End - 1)
$spaceRange.Delete()

# Update the placeholder text to reference the new section id.
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("**ID__AFFARS_pgi_5339_topic_3__ID**", $true, $false, $false, $false, $false,
                        $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5339_101_90__ID**", 2)

# Re-fetch the paragraph/format handles after the text edits above.
$p1 = $d.Paragraphs(1)
$fmt = $p1.Format

# w:ind w:left="225" (twips) == 11.25pt
$fmt.LeftIndent = 11.25

# Add a w:pBdr with w:space="5" on every side (no line drawn).
$borders = $fmt.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
